# Daily attendance processing - 2026-02-07 05:35:04 UTC
# Swap the "Recorded By" name ordering in column G from
# "Administrator, Miss Dina Nasr" to "Miss Dina Nasr, Administrator"
# wherever it appears in the workbook.

$wb = $excel.ActiveWorkbook
$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            if ($cell.Value2 -eq $oldText) {
                $cell.Value2 = $newText
            }
        }
    }
}
